# Auto-generated PowerShell-style Excel COM-interop script
# Applies profit-table recalculation updates across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1909
$ws.Range("J29").Value = 2613.8
$ws.Range("L29").Value = 7841.400000000001
$ws.Range("N29").Value = -8403.400000000001
$ws.Range("H38").Value = 253
$ws.Range("I38").Value = 159.625
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 478.875
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -106.875
$ws.Range("N38").Value = -3744
$ws.Range("H39").Value = 666.2273
$ws.Range("I39").Value = 747.1818
$ws.Range("K39").Value = 2241.5454
$ws.Range("M39").Value = -1945.5454
$ws.Range("H58").Value = 3295
$ws.Range("J58").Value = 5250
$ws.Range("L58").Value = 15750
$ws.Range("N58").Value = -16050
$ws.Range("H132").Value = 8452.9375
$ws.Range("I132").Value = 4824.7
$ws.Range("K132").Value = 14474.1
$ws.Range("M132").Value = -11944.1
$ws.Range("H137").Value = 1275.7407
$ws.Range("I137").Value = 1062.6154
$ws.Range("J137").Value = 1473.6428
$ws.Range("K137").Value = 3187.8462
$ws.Range("L137").Value = 4420.928400000001
$ws.Range("M137").Value = -637.8462
$ws.Range("N137").Value = -9520.928400000001
$ws.Range("H138").Value = 2022.3334
$ws.Range("J138").Value = 2103.6628
$ws.Range("L138").Value = 6310.9884
$ws.Range("N138").Value = -16590.9884

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3782.1064
$ws.Range("I32").Value = 3782.1064
$ws.Range("K32").Value = 3782.1064
$ws.Range("M32").Value = -3495.1064
$ws.Range("H45").Value = 1693.3684
$ws.Range("I45").Value = 1598.4706
$ws.Range("K45").Value = 1598.4706
$ws.Range("M45").Value = -1221.4706
$ws.Range("H61").Value = 1296
$ws.Range("I61").Value = 1296
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1296
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1084
$ws.Range("N61").ClearContents()
$ws.Range("H97").Value = 457.77777
$ws.Range("I97").Value = 457.77777
$ws.Range("K97").Value = 457.77777
$ws.Range("M97").Value = 38.22223000000002
$ws.Range("H132").Value = 3919.9375
$ws.Range("I132").Value = 4322.2
$ws.Range("J132").Value = 3249.5
$ws.Range("K132").Value = 12966.6
$ws.Range("L132").Value = 9748.5
$ws.Range("M132").Value = -10436.6
$ws.Range("N132").Value = -14808.5
$ws.Range("H136").Value = 1296
$ws.Range("I136").Value = 1296
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3888
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1338
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4554.067
$ws.Range("I86").Value = 4864.636
$ws.Range("K86").Value = 4864.636
$ws.Range("M86").Value = -3741.636
$ws.Range("H89").Value = 4554.067
$ws.Range("I89").Value = 4864.636
$ws.Range("K89").Value = 24323.18
$ws.Range("M89").Value = -18707.18
$ws.Range("H94").Value = 31250894
$ws.Range("I94").Value = 50000532
$ws.Range("J94").Value = 1499.6666
$ws.Range("K94").Value = 50000532
$ws.Range("L94").Value = 1499.6666
$ws.Range("M94").Value = -50000081
$ws.Range("N94").Value = -2401.6666
$ws.Range("H134").Value = 4919.2964
$ws.Range("I134").Value = 1061.3043
$ws.Range("K134").Value = 3183.9129
$ws.Range("M134").Value = -648.9129000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1717.4762
$ws.Range("J31").Value = 2028.4615
$ws.Range("L31").Value = 2028.4615
$ws.Range("N31").Value = -2618.4615
$ws.Range("H34").Value = 1717.4762
$ws.Range("J34").Value = 2028.4615
$ws.Range("L34").Value = 2028.4615
$ws.Range("N34").Value = -2432.4615
$ws.Range("H39").Value = 999
$ws.Range("I39").Value = 999
$ws.Range("K39").Value = 999
$ws.Range("M39").Value = -608
$ws.Range("H49").Value = 999
$ws.Range("I49").Value = 999
$ws.Range("K49").Value = 999
$ws.Range("M49").Value = -817
$ws.Range("H132").Value = 7106.6665
$ws.Range("I132").Value = 9217.538
$ws.Range("K132").Value = 27652.614
$ws.Range("M132").Value = -25122.614
$ws.Range("H134").Value = 12347255
$ws.Range("I134").Value = 15874428
$ws.Range("J134").Value = 2149.6667
$ws.Range("K134").Value = 47623284
$ws.Range("L134").Value = 6449.000100000001
$ws.Range("M134").Value = -47620749
$ws.Range("N134").Value = -11519.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3000.5
$ws.Range("J55").Value = 3285.7144
$ws.Range("L55").Value = 9857.143199999999
$ws.Range("N55").Value = -10211.1432
$ws.Range("H68").Value = 1422.909
$ws.Range("I68").Value = 764.44446
$ws.Range("J68").Value = 1669.8334
$ws.Range("K68").Value = 2293.33338
$ws.Range("L68").Value = 5009.5002
$ws.Range("M68").Value = -1482.33338
$ws.Range("N68").Value = -6631.5002
$ws.Range("H71").Value = 1422.909
$ws.Range("I71").Value = 764.44446
$ws.Range("J71").Value = 1669.8334
$ws.Range("K71").Value = 6880.00014
$ws.Range("L71").Value = 15028.5006
$ws.Range("M71").Value = -2824.00014
$ws.Range("N71").Value = -23140.5006
$ws.Range("H107").Value = 5686.5
$ws.Range("J107").Value = 9704.454
$ws.Range("L107").Value = 29113.362
$ws.Range("N107").Value = -32953.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 21000
$ws.Range("J46").Value = 21000
$ws.Range("L46").Value = 21000
$ws.Range("N46").Value = -21312
$ws.Range("H57").Value = 24999.8
$ws.Range("J57").Value = 24999.8
$ws.Range("L57").Value = 24999.8
$ws.Range("N57").Value = -26639.8
$ws.Range("H113").Value = 1433.0454
$ws.Range("I113").Value = 901.58826
$ws.Range("J113").Value = 3240
$ws.Range("K113").Value = 901.58826
$ws.Range("L113").Value = 3240
$ws.Range("M113").Value = 1268.41174
$ws.Range("N113").Value = -7580
$ws.Range("H132").Value = 1944.3334
$ws.Range("I132").Value = 1636.5
$ws.Range("J132").Value = 2560
$ws.Range("K132").Value = 4909.5
$ws.Range("L132").Value = 7680
$ws.Range("M132").Value = -2379.5
$ws.Range("N132").Value = -12740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 768
$ws.Range("I22").Value = 344.125
$ws.Range("J22").Value = 1333.1666
$ws.Range("K22").Value = 344.125
$ws.Range("L22").Value = 1333.1666
$ws.Range("M22").Value = -49.125
$ws.Range("N22").Value = -1923.1666
$ws.Range("H27").Value = 768
$ws.Range("I27").Value = 344.125
$ws.Range("J27").Value = 1333.1666
$ws.Range("K27").Value = 344.125
$ws.Range("L27").Value = 1333.1666
$ws.Range("M27").Value = -237.125
$ws.Range("N27").Value = -1547.1666
$ws.Range("H132").Value = 54911.367
$ws.Range("I132").Value = 1816.3846
$ws.Range("J132").Value = 169950.5
$ws.Range("K132").Value = 5449.1538
$ws.Range("L132").Value = 509851.5
$ws.Range("M132").Value = -2919.1538
$ws.Range("N132").Value = -514911.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H126").Value = 62501460
$ws.Range("I126").Value = 111111810
$ws.Range("K126").Value = 333335430
$ws.Range("M126").Value = -333332960
$ws.Range("H132").Value = 3679.2856
$ws.Range("I132").Value = 4467.8335
$ws.Range("K132").Value = 13403.5005
$ws.Range("M132").Value = -10873.5005
$ws.Range("H136").Value = 700.13043
$ws.Range("I136").Value = 439.95
$ws.Range("J136").Value = 2434.6667
$ws.Range("K136").Value = 2434.6667
$ws.Range("L136").Value = 7304.000100000001
$ws.Range("M136").Value = 1230.15
$ws.Range("N136").Value = -12404.0001
